$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 14 new rows before the old row 65 (pushes old row 65 -> row 79)
$ws.Rows("65:78").Insert()

# Row 64: new block title
$ws.Range("A64").Value = "Multiresolution Histograms (3nd Run - no SVM, bayes NB_THRESH = .5, db4)"
$ws.Range("A64").Font.Bold = $true

# Row 65: "Guess" label
$ws.Range("B65").Value = "Guess"
$ws.Range("B65").Font.Bold = $true

# Row 66: header row (category names + Percentage/False Neg/True Pos)
$ws.Range("B66").Value = "Airport"
$ws.Range("C66").Value = "Auditorium"
$ws.Range("D66").Value = "Bamboo"
$ws.Range("E66").Value = "Campus"
$ws.Range("F66").Value = "Desert"
$ws.Range("G66").Value = "Football Field"
$ws.Range("H66").Value = "Kitchen"
$ws.Range("I66").Value = "Sky"
$ws.Range("J66").Value = "Percentage"
$ws.Range("K66").Value = "False Neg"
$ws.Range("L66").Value = "True Pos"

# Rows 67-74: confusion matrix data + J/K/L formulas
$ws.Range("A67").Value = "Airport"
$ws.Range("B67").Value = 1
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 0
$ws.Range("E67").Value = 2
$ws.Range("F67").Value = 13
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 1
$ws.Range("I67").Value = 2
$ws.Range("K67").Formula = "=(SUM(B67:I67) - B67) / SUM(B67:I67)"

$ws.Range("A68").Value = "Auditorium"
$ws.Range("B68").Value = 2
$ws.Range("C68").Value = 10
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 4
$ws.Range("I68").Value = 1
$ws.Range("K68").Formula = "=(SUM(B68:I68) - C68) / SUM(B68:I68)"

$ws.Range("A69").Value = "Bamboo"
$ws.Range("B69").Value = 0
$ws.Range("C69").Value = 1
$ws.Range("D69").Value = 15
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 3
$ws.Range("I69").Value = 1
$ws.Range("K69").Formula = "=(SUM(B69:I69) - D69) / SUM(B69:I69)"

$ws.Range("A70").Value = "Campus"
$ws.Range("B70").Value = 0
$ws.Range("C70").Value = 8
$ws.Range("D70").Value = 10
$ws.Range("E70").Value = 1
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 1
$ws.Range("K70").Formula = "=(SUM(B70:I70) - E70) / SUM(B70:I70)"

$ws.Range("A71").Value = "Desert"
$ws.Range("B71").Value = 0
$ws.Range("C71").Value = 2
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 1
$ws.Range("F71").Value = 16
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 1
$ws.Range("I71").Value = 0
$ws.Range("K71").Formula = "=(SUM(B71:I71) - F71) / SUM(B71:I71)"

$ws.Range("A72").Value = "Football Field"
$ws.Range("B72").Value = 3
$ws.Range("C72").Value = 6
$ws.Range("D72").Value = 3
$ws.Range("E72").Value = 0
$ws.Range("F72").Value = 3
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 2
$ws.Range("I72").Value = 1
$ws.Range("K72").Formula = "=(SUM(B72:I72) - G72) / SUM(B72:I72)"

$ws.Range("A73").Value = "Kitchen"
$ws.Range("B73").Value = 1
$ws.Range("C73").Value = 4
$ws.Range("D73").Value = 1
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 12
$ws.Range("I73").Value = 0
$ws.Range("K73").Formula = "=(SUM(B73:I73) - H73) / SUM(B73:I73)"

$ws.Range("A74").Value = "Sky"
$ws.Range("B74").Value = 0
$ws.Range("C74").Value = 2
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 6
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 1
$ws.Range("I74").Value = 11
$ws.Range("K74").Formula = "=(SUM(B74:I74) - I74) / SUM(B74:I74)"

# Column J (rows 67:74) as one fill (Excel will create shared formula J68:J74, J67 separate)
$ws.Range("J67").Formula = "=SUM(B67:I67)/SUM(`$B`$22:`$I`$29)"
$ws.Range("J68:J74").Formula = "=SUM(B68:I68)/SUM(`$B`$22:`$I`$29)"

# Column L (rows 67:74)
$ws.Range("L67").Formula = "=1-K67"
$ws.Range("L68:L74").Formula = "=1-K68"

# Row 75: percentage totals
$ws.Range("A75").Value = "Percentage"
$ws.Range("B75").Formula = "=SUM(B67:B74) / SUM(`$B`$52:`$I`$59)"
$ws.Range("C75:I75").Formula = "=SUM(C67:C74) / SUM(`$B`$52:`$I`$59)"

# Row 76: accuracy per category + K76 label + L76 total accuracy
$ws.Range("A76").Value = "False Pos"
$ws.Range("B76").Formula = "=(SUM(B67:B74) - B67) / SUM(B67:B74)"
$ws.Range("C76").Formula = "=(SUM(C67:C74) - C68) / SUM(C67:C74)"
$ws.Range("D76").Formula = "=(SUM(D67:D74) - D69) / SUM(D67:D74)"
$ws.Range("E76").Formula = "=(SUM(E67:E74) - E69) / SUM(E67:E74)"
$ws.Range("F76").Formula = "=(SUM(F67:F74) - F71) / SUM(F67:F74)"
$ws.Range("G76").Formula = "=(SUM(G67:G74) - G72) / SUM(G67:G74)"
$ws.Range("H76").Formula = "=(SUM(H67:H74) - H73) / SUM(H67:H74)"
$ws.Range("I76").Formula = "=(SUM(I67:I74) - I74) / SUM(I67:I74)"
$ws.Range("K76").Value = "Accuracy"
$ws.Range("K76").Font.Bold = $true
$ws.Range("L76").Formula = "=(B67+C68+D69+E70+F71+G72+H73+I74) / SUM(B67:I74)"

# Row 77: error rate per category
$ws.Range("A77").Value = "True Neg"
$ws.Range("B77").Formula = "=1-B76"
$ws.Range("C77").Formula = "=1-C76"
$ws.Range("D77:I77").Formula = "=1-D76"

# Row 79: trailing title (moved from old row 65), restore bold
$ws.Range("A79").Font.Bold = $true

# Update sheet view scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 60
$ws.Range("E60").Select()
